$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.181.31'
$ws.Range("E2").Value = '  +1.57%  '

$ws.Range("D3").Value = '1.877.01'
$ws.Range("E3").Value = '  +1.49%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.21'
$ws.Range("E5").Value = '  +0.31%  '

$ws.Range("E6").Value = '  -0.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4313'
$ws.Range("E7").Value = '  +1.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3703'
$ws.Range("E8").Value = '  +1.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07438'
$ws.Range("E9").Value = '  +1.57%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8957'
$ws.Range("E10").Value = '  +0.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.17'
$ws.Range("E11").Value = '  +1.41%  '

$ws.Range("D12").Value = '1.969.50'
$ws.Range("E12").Value = '  +3.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.444'
$ws.Range("E13").Value = '  +1.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.663'
$ws.Range("E14").Value = '  +1.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06995'
$ws.Range("E15").Value = '  +0.88%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").Value = '  -0.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '81.41'
$ws.Range("E17").Value = '  +2.78%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009104'
$ws.Range("E18").Value = '  +2.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.67'
$ws.Range("E20").Value = '  +1.06%  '

$ws.Range("D21").Value = '28.149.26'
$ws.Range("E21").Value = '  +1.45%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.086'
$ws.Range("E22").Value = '  +1.89%  '

$ws.Range("E23").Value = '  +2.58%  '

$ws.Range("D24").Value = '2.140.66'
$ws.Range("E24").Value = '  +2.65%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.988'
$ws.Range("E25").Value = '  +1.25%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.20'
$ws.Range("E26").Value = '  +0.44%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.79'
$ws.Range("E27").Value = '  -0.90%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.416'
$ws.Range("E28").Value = '  +3.19%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '118.49'
$ws.Range("E29").Value = '  -2.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.912'
$ws.Range("E30").Value = '  +0.35%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08946'
$ws.Range("E31").Value = '  +0.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7988'
$ws.Range("E32").Value = '  +3.71%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.687'
$ws.Range("E33").Value = '  +2.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.182'
$ws.Range("E34").Value = '  +7.02%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.971'
$ws.Range("E35").Value = '  +0.04%  '

$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.130'
$ws.Range("E36").Value = '  +3.10%  '

$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.001'
$ws.Range("E37").Value = '  -0.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05499'
$ws.Range("E38").Value = '  +2.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01966'
$ws.Range("E39").Value = '  +0.44%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.900'
$ws.Range("E40").Value = '  +3.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5191'
$ws.Range("E41").Value = '  +1.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1697'
$ws.Range("E42").Value = '  +1.99%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.867'
$ws.Range("E43").Value = '  -0.74%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.589'
$ws.Range("E44").Value = '  +3.65%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.59'
$ws.Range("E45").Value = '  +1.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06605'
$ws.Range("E46").Value = '  +0.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4779'
$ws.Range("E47").Value = '  +0.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '106.01'
$ws.Range("E48").Value = '  +1.19%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.000'
$ws.Range("E49").Value = '  -0.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.661'
$ws.Range("E50").Value = '  +1.52%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.889'
$ws.Range("E51").Value = '  +6.66%  '
